# Edit slide 1 "TextBox 5" quote content + reposition/resize the box
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$shp = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $candidate = $s.Shapes.Item($i)
    if ($candidate.Name -eq "TextBox 5") {
        $shp = $candidate
    }
}

$tr = $shp.TextFrame.TextRange
$tr.Text = "“We call these algorithms data parallel algorithms because their parallelism comes from simultaneous operations across large sets of data, rather than from multiple thread of control.”`r`t- W. Daniel Hillis and Guy L. Steele`r`t“Data Parallel Algorithms,” Comm. ACM (1986)`r`r“If you were plowing a field, which would you rather use, two strong oxen or 1024 chickens?”`r`t- Seymour Cray, Father of the Supercomputer`r`t(arguing for two powerful vector processors `r`tversus many simple processors)"

# Italicize "data parallel " inside paragraph 1 and "Comm. ACM " inside paragraph 3
$tr.Characters(27, 14).Font.Italic = $true
$tr.Characters(253, 10).Font.Italic = $true

# Reposition / resize the textbox (values nudged so the Single-precision
# round trip lands exactly on the target EMU values)
$shp.Left = 281.0025196850394
$shp.Top = 320.68583877165355
$shp.Width = 404.12504737007873
$shp.Height = 193.87503937007875
